$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '51.758.92'
Set-TextValue 'E2' '  +0.50%  '
Set-TextValue 'D3' '2.832.43'
Set-TextValue 'E3' '  +2.67%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '353.51'
Set-TextValue 'E5' '  +6.09%  '
Set-TextValue 'D6' '113.68'
Set-TextValue 'E6' '  -2.18%  '
Set-TextValue 'D7' '0.572'
Set-TextValue 'E7' '  +6.25%  '
Set-TextValue 'D8' '0.999'
Set-TextValue 'E8' '  -0.02%  '
Set-TextValue 'D9' '0.599'
Set-TextValue 'E9' '  +4.43%  '
Set-TextValue 'D10' '41.68'
Set-TextValue 'E10' '  -0.16%  '
Set-TextValue 'D11' '0.0855'
Set-TextValue 'E11' '  -1.54%  '
Set-TextValue 'D12' '20.06'
Set-TextValue 'E12' '  -1.17%  '
Set-TextValue 'E13' '  +1.34%  '
Set-TextValue 'E14' '  +0.82%  '
Set-TextValue 'D15' '3.279.45'
Set-TextValue 'E15' '  +2.54%  '
Set-TextValue 'D16' '2.830.46'
Set-TextValue 'E16' '  +2.54%  '
Set-TextValue 'D17' '0.899'
Set-TextValue 'E17' '  +1.11%  '
Set-TextValue 'D18' '51.614.42'
Set-TextValue 'E18' '  +0.13%  '
Set-TextValue 'D19' '7.36'
Set-TextValue 'E19' '  +7.19%  '
Set-TextValue 'D20' '3.17'
Set-TextValue 'E20' '  -3.30%  '
Set-TextValue 'D21' '13.51'
Set-TextValue 'E21' '  +0.19%  '
Set-TextValue 'D22' '0.0₃0997'
Set-TextValue 'E22' '  +2.02%  '
Set-TextValue 'D23' '271.03'
Set-TextValue 'E23' '  -2.69%  '
Set-TextValue 'D24' '69.80'
Set-TextValue 'E24' '  +0.18%  '
Set-TextValue 'D25' '2.77'
Set-TextValue 'E25' '  +3.41%  '
Set-TextValue 'D26' '26.76'
Set-TextValue 'E26' '  -0.13%  '
Set-TextValue 'E27' '  +0.01%  '
Set-TextValue 'E28' '  +1.50%  '
Set-TextValue 'E29' '  +1.10%  '
Set-TextValue 'E30' '  -2.07%  '
Set-TextValue 'D31' '34.22'
Set-TextValue 'E31' '  -2.43%  '
Set-TextValue 'D32' '50.69'
Set-TextValue 'E32' '  +1.23%  '
Set-TextValue 'E33' '  +4.78%  '
Set-TextValue 'D34' '0.0442'
Set-TextValue 'E34' '  +25.44%  '
Set-TextValue 'D35' '0.0825'
Set-TextValue 'E35' '  +0.21%  '
Set-TextValue 'D36' '1.00'
Set-TextValue 'E36' '  -0.17%  '
Set-TextValue 'D37' '2.07'
Set-TextValue 'E37' '  -0.45%  '
Set-TextValue 'D38' '4.90'
Set-TextValue 'E38' '  -1.96%  '
Set-TextValue 'D39' '3.19'
Set-TextValue 'E39' '  -1.51%  '
Set-TextValue 'D40' '18.13'
Set-TextValue 'E40' '  -4.45%  '
Set-TextValue 'D41' '23.81'
Set-TextValue 'E41' '  +3.25%  '
Set-TextValue 'D42' '0.117'
Set-TextValue 'E42' '  +2.44%  '
Set-TextValue 'D43' '2.53'
Set-TextValue 'E43' '  +3.13%  '
Set-TextValue 'D44' '126.12'
Set-TextValue 'E44' '  -1.00%  '
Set-TextValue 'E45' '  -0.44%  '
Set-TextValue 'D46' '2.080.01'
Set-TextValue 'E46' '  -0.52%  '
Set-TextValue 'E47' '  +0.44%  '
Set-TextValue 'E48' '  +3.68%  '
Set-TextValue 'E49' '  +2.93%  '
Set-TextValue 'D50' '0.934'
Set-TextValue 'E50' '  +6.82%  '
Set-TextValue 'D51' '61.03'
Set-TextValue 'E51' '  +1.55%  '
